$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.410.02'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '1.728.57'
$ws.Range("E3").Value = '  +3.34%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").Value = '''243.10'
$ws.Range("E5").Value = '  +2.35%  '
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").Value = '''0.4787'
$ws.Range("E7").Value = '  +3.50%  '
$ws.Range("D8").Value = '''0.2660'
$ws.Range("E8").Value = '  +2.93%  '
$ws.Range("D9").Value = '''0.06224'
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("D10").Value = '1.731.08'
$ws.Range("E10").Value = '  +4.00%  '
$ws.Range("D11").Value = '''0.07112'
$ws.Range("E11").Value = '  +2.48%  '
$ws.Range("D12").Value = '''15.66'
$ws.Range("E12").Value = '  +4.66%  '
$ws.Range("D13").Value = '''0.6153'
$ws.Range("E13").Value = '  +6.70%  '
$ws.Range("D14").Value = '''4.560'
$ws.Range("E14").Value = '  +4.63%  '
$ws.Range("D15").Value = '''76.79'
$ws.Range("E15").Value = '  +2.01%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '26.435.56'
$ws.Range("E17").Value = '  +2.86%  '
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '''0.000006898'
$ws.Range("E19").Value = '  +2.78%  '
$ws.Range("E20").Value = '  +2.61%  '
$ws.Range("D21").Value = '1.957.49'
$ws.Range("E21").Value = '  +4.27%  '
$ws.Range("D22").Value = '''4.568'
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = '''8.876'
$ws.Range("E23").Value = '  +2.79%  '
$ws.Range("D24").Value = '''5.312'
$ws.Range("E24").Value = '  +1.33%  '
$ws.Range("D25").Value = '''136.13'
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").Value = '''1.786'
$ws.Range("E27").Value = '  +3.45%  '
$ws.Range("D28").Value = '''1.404'
$ws.Range("E28").Value = '  -0.95%  '
$ws.Range("D29").Value = '''106.36'
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("D30").Value = '''3.972'
$ws.Range("E30").Value = '  +0.50%  '
$ws.Range("D31").Value = '''0.07988'
$ws.Range("E31").Value = '  +4.27%  '
$ws.Range("D32").Value = '''3.711'
$ws.Range("E32").Value = '  +2.88%  '
$ws.Range("D33").Value = '''0.04535'
$ws.Range("E33").Value = '  +3.98%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.616'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '''0.6364'
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''0.9851'
$ws.Range("E36").Value = '  +4.25%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''0.9318'
$ws.Range("E37").Value = '  +0.77%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").Value = '''109.04'
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").Value = '''1.975'
$ws.Range("E39").Value = '  +7.51%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.403'
$ws.Range("E40").Value = '  -1.71%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").Value = '''1.006'
$ws.Range("E41").Value = '  +0.79%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01504'
$ws.Range("E42").Value = '  +3.79%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '''5.630'
$ws.Range("E43").Value = '  +11.36%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.3901'
$ws.Range("E44").Value = '  +4.89%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '''6.922'
$ws.Range("E45").Value = '  +13.06%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '''0.1191'
$ws.Range("E46").Value = '  +6.80%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.05333'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '''30.77'
$ws.Range("E48").Value = '  -1.09%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''7.844'
$ws.Range("E49").Value = '  +2.83%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.273'
$ws.Range("E50").Value = '  +5.37%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = '''0.3419'
$ws.Range("E51").Value = '  +2.83%  '
